# Scheduled-runner profit recalculation update for Sheets/Mateus_Profits.xlsx
# Updates currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ /
# LeveProfitNQ/HQ figures across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and
# WVR market-data tables to reflect refreshed pricing.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3847033.5
$ws.Range("J17").Value = 3847033.5
$ws.Range("L17").Value = 11541100.5
$ws.Range("N17").Value = -11541436.5

$ws.Range("H62").Value = 26166.834
$ws.Range("I62").Value = 26500.25
$ws.Range("K62").Value = 26500.25
$ws.Range("M62").Value = -25876.25

$ws.Range("H65").Value = 26166.834
$ws.Range("I65").Value = 26500.25
$ws.Range("K65").Value = 132501.25
$ws.Range("M65").Value = -129381.25

$ws.Range("H132").Value = 1785.52
$ws.Range("I132").Value = 1739.6818
$ws.Range("K132").Value = 5219.0454
$ws.Range("M132").Value = -2689.0454

$ws.Range("H138").Value = 23258754
$ws.Range("I138").Value = 4396.75
$ws.Range("J138").Value = 25643816
$ws.Range("K138").Value = 13190.25
$ws.Range("L138").Value = 76931448
$ws.Range("M138").Value = -8050.25
$ws.Range("N138").Value = -76941728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5096.074
$ws.Range("I32").Value = 5000.451
$ws.Range("K32").Value = 5000.451
$ws.Range("M32").Value = -4713.451

$ws.Range("H61").Value = 22734300
$ws.Range("I61").Value = 26321558
$ws.Range("K61").Value = 26321558
$ws.Range("M61").Value = -26321346

$ws.Range("H74").Value = 3807
$ws.Range("I74").Value = 3228.0908
$ws.Range("K74").Value = 3228.0908
$ws.Range("M74").Value = -2354.0908

$ws.Range("H77").Value = 3807
$ws.Range("I77").Value = 3228.0908
$ws.Range("K77").Value = 16140.454
$ws.Range("M77").Value = -11772.454

$ws.Range("H122").Value = 2487.2222
$ws.Range("J122").Value = 1995
$ws.Range("L122").Value = 5985
$ws.Range("N122").Value = -10885

$ws.Range("H136").Value = 22734300
$ws.Range("I136").Value = 26321558
$ws.Range("K136").Value = 78964674
$ws.Range("M136").Value = -78962124

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3983.6667
$ws.Range("J105").Value = 4223
$ws.Range("L105").Value = 4223
$ws.Range("N105").Value = -7717

$ws.Range("H134").Value = 4838.4707
$ws.Range("I134").Value = 4515.875
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 13547.625
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -11012.625
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9997.75
$ws.Range("I25").Value = 9997.75
$ws.Range("K25").Value = 9997.75
$ws.Range("M25").Value = -9823.75

$ws.Range("H31").Value = 6212.8237
$ws.Range("I31").Value = 4936.5835
$ws.Range("K31").Value = 4936.5835
$ws.Range("M31").Value = -4641.5835

$ws.Range("H34").Value = 6212.8237
$ws.Range("I34").Value = 4936.5835
$ws.Range("K34").Value = 4936.5835
$ws.Range("M34").Value = -4734.5835

$ws.Range("H58").Value = 11999.889
$ws.Range("I58").Value = 7500
$ws.Range("K58").Value = 7500
$ws.Range("M58").Value = -7297

$ws.Range("H96").Value = 52437.5
$ws.Range("J96").Value = 52437.5
$ws.Range("L96").Value = 52437.5
$ws.Range("N96").Value = -57929.5

$ws.Range("H134").Value = 6289.2666
$ws.Range("I134").Value = 5278.25
$ws.Range("K134").Value = 15834.75
$ws.Range("M134").Value = -13299.75

$ws.Range("H136").Value = 11999.889
$ws.Range("I136").Value = 7500
$ws.Range("K136").Value = 22500
$ws.Range("M136").Value = -19950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 800
$ws.Range("I31").Value = 800
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2400
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2112
$ws.Range("N31").ClearContents()

$ws.Range("H112").Value = 5000
$ws.Range("J112").Value = 5000
$ws.Range("L112").Value = 15000
$ws.Range("N112").Value = -17216

$ws.Range("H113").Value = 489.8889
$ws.Range("J113").Value = 601.36365
$ws.Range("L113").Value = 1804.09095
$ws.Range("N113").Value = -6144.09095

$ws.Range("H122").Value = 737.8182
$ws.Range("I122").Value = 388.66666
$ws.Range("K122").Value = 3497.99994
$ws.Range("M122").Value = -1047.99994

$ws.Range("H129").Value = 10645.143
$ws.Range("I129").Value = 2308
$ws.Range("K129").Value = 6924
$ws.Range("M129").Value = -1924

$ws.Range("H131").Value = 27780278
$ws.Range("I131").Value = 166667330
$ws.Range("K131").Value = 500001990
$ws.Range("M131").Value = -499996950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4546.4
$ws.Range("I5").Value = 4433
$ws.Range("K5").Value = 4433
$ws.Range("M5").Value = -4321

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H47").Value = 24999.5
$ws.Range("J47").Value = 24999.5
$ws.Range("L47").Value = 24999.5
$ws.Range("N47").Value = -26135.5

$ws.Range("H102").Value = 3257.389
$ws.Range("J102").Value = 6484.5
$ws.Range("L102").Value = 6484.5
$ws.Range("N102").Value = -9728.5

$ws.Range("H122").Value = 2270.7144
$ws.Range("I122").Value = 2270.7144
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6812.1432
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4362.1432
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4399.5
$ws.Range("J40").Value = 4999.5
$ws.Range("L40").Value = 4999.5
$ws.Range("N40").Value = -5271.5

$ws.Range("H82").Value = 15366.714
$ws.Range("I82").Value = 17113.4
$ws.Range("J82").Value = 11000
$ws.Range("K82").Value = 17113.4
$ws.Range("L82").Value = 11000
$ws.Range("M82").Value = -16752.4
$ws.Range("N82").Value = -11722

$ws.Range("H85").Value = 15366.714
$ws.Range("I85").Value = 17113.4
$ws.Range("J85").Value = 11000
$ws.Range("K85").Value = 17113.4
$ws.Range("L85").Value = 11000
$ws.Range("M85").Value = -15865.4
$ws.Range("N85").Value = -13496

$ws.Range("H97").Value = 24000
$ws.Range("J97").Value = 24000
$ws.Range("L97").Value = 24000
$ws.Range("N97").Value = -25982

$ws.Range("H100").Value = 2179347.8
$ws.Range("I100").Value = 3129743.5
$ws.Range("J100").Value = 7014.143
$ws.Range("K100").Value = 3129743.5
$ws.Range("L100").Value = 7014.143
$ws.Range("M100").Value = -3129202.5
$ws.Range("N100").Value = -8096.143

$ws.Range("H122").Value = 2627.6667
$ws.Range("I122").Value = 3015.8
$ws.Range("K122").Value = 9047.400000000001
$ws.Range("M122").Value = -6597.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 222222580
$ws.Range("I6").Value = 222222580
$ws.Range("K6").Value = 222222580
$ws.Range("M6").Value = -222222465

$ws.Range("H81").Value = 2521.818
$ws.Range("I81").Value = 2012.9333
$ws.Range("K81").Value = 4025.8666
$ws.Range("M81").Value = -2964.8666

$ws.Range("H84").Value = 2521.818
$ws.Range("I84").Value = 2012.9333
$ws.Range("K84").Value = 20129.333
$ws.Range("M84").Value = -14825.333

$ws.Range("H122").Value = 3841.476
$ws.Range("I122").Value = 2489.3333
$ws.Range("K122").Value = 7467.999899999999
$ws.Range("M122").Value = -5017.999899999999

$ws.Range("H132").Value = 4316.6816
$ws.Range("I132").Value = 4046.0476
$ws.Range("K132").Value = 12138.1428
$ws.Range("M132").Value = -9608.1428
